$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 82 (G=12623)
$ws.Range("H82").Value = 947.3333
$ws.Range("I82").Value = 671
$ws.Range("K82").Value = 2013
$ws.Range("M82").Value = -1607

# Row 85 (G=12623)
$ws.Range("H85").Value = 947.3333
$ws.Range("I85").Value = 671
$ws.Range("K85").Value = 2013
$ws.Range("M85").Value = -609

# Row 94 (G=19905)
$ws.Range("H94").Value = 531.7143
$ws.Range("I94").Value = 531.7143
$ws.Range("K94").Value = 531.7143
$ws.Range("M94").Value = -80.71429999999998

# Row 98 (G=36237)
$ws.Range("H98").Value = 2221.6
$ws.Range("I98").Value = 1736
$ws.Range("J98").Value = 2950
$ws.Range("K98").Value = 1736
$ws.Range("L98").Value = 2950
$ws.Range("M98").Value = -238
$ws.Range("N98").Value = -5946

# Row 107 (G=27766)
$ws.Range("H107").Value = 2683.375
$ws.Range("I107").Value = 2175.5
$ws.Range("J107").Value = 5222.75
$ws.Range("K107").Value = 2175.5
$ws.Range("L107").Value = 5222.75
$ws.Range("M107").Value = -255.5
$ws.Range("N107").Value = -9062.75

# Row 122 (G=36237)
$ws.Range("H122").Value = 2221.6
$ws.Range("I122").Value = 1736
$ws.Range("J122").Value = 2950
$ws.Range("K122").Value = 5208
$ws.Range("L122").Value = 8850
$ws.Range("M122").Value = -2758
$ws.Range("N122").Value = -13750

# Row 131 (G=36108)
$ws.Range("H131").Value = 5465.154
$ws.Range("I131").Value = 1753.9166
$ws.Range("K131").Value = 5261.7498
$ws.Range("M131").Value = -221.7497999999996

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (G=44147)
$ws.Range("H32").Value = 6159.9473
$ws.Range("I32").Value = 5047.4517
$ws.Range("K32").Value = 5047.4517
$ws.Range("M32").Value = -4760.4517

# Row 137 (G=43227)
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 58 (G=43234)
$ws.Range("H58").Value = 38999.5
$ws.Range("J58").Value = 38000
$ws.Range("L58").Value = 38000
$ws.Range("N58").Value = -38588

# Row 105 (G=19947)
$ws.Range("H105").Value = 3238.8572
$ws.Range("I105").Value = 3475.375
$ws.Range("J105").Value = 2482
$ws.Range("K105").Value = 3475.375
$ws.Range("L105").Value = 2482
$ws.Range("M105").Value = -1728.375
$ws.Range("N105").Value = -5976

# Row 115 (G=27118)
$ws.Range("H115").Value = 50000
$ws.Range("J115").Value = 50000
$ws.Range("L115").Value = 50000
$ws.Range("N115").Value = -53134

# Row 134 (G=43998)
$ws.Range("H134").Value = 9194.821
$ws.Range("I134").Value = 5920.3335
$ws.Range("J134").Value = 15088.9
$ws.Range("K134").Value = 17761.0005
$ws.Range("L134").Value = 45266.7
$ws.Range("M134").Value = -15226.0005
$ws.Range("N134").Value = -50336.7

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (G=44023)
$ws.Range("H31").Value = 2089.1667
$ws.Range("I31").Value = 1347.5
$ws.Range("J31").Value = 2460
$ws.Range("K31").Value = 1347.5
$ws.Range("L31").Value = 2460
$ws.Range("M31").Value = -1052.5
$ws.Range("N31").Value = -3050

# Row 34 (G=44023)
$ws.Range("H34").Value = 2089.1667
$ws.Range("I34").Value = 1347.5
$ws.Range("J34").Value = 2460
$ws.Range("K34").Value = 1347.5
$ws.Range("L34").Value = 2460
$ws.Range("M34").Value = -1145.5
$ws.Range("N34").Value = -2864

# Row 99 (G=36198)
$ws.Range("H99").Value = 9154.166999999999
$ws.Range("I99").Value = 6017.625
$ws.Range("J99").Value = 9892.177
$ws.Range("K99").Value = 6017.625
$ws.Range("L99").Value = 9892.177
$ws.Range("M99").Value = -4519.625
$ws.Range("N99").Value = -12888.177

# Row 107 (G=27689)
$ws.Range("H107").Value = 881.3125
$ws.Range("I107").Value = 687.3684
$ws.Range("K107").Value = 687.3684
$ws.Range("M107").Value = 1232.6316

# Row 126 (G=36198)
$ws.Range("H126").Value = 9154.166999999999
$ws.Range("I126").Value = 6017.625
$ws.Range("J126").Value = 9892.177
$ws.Range("K126").Value = 18052.875
$ws.Range("L126").Value = 29676.531
$ws.Range("M126").Value = -15582.875
$ws.Range("N126").Value = -34616.531

# Row 132 (G=44019)
$ws.Range("H132").Value = 26571.441
$ws.Range("J132").Value = 41100.285
$ws.Range("L132").Value = 123300.855
$ws.Range("N132").Value = -128360.855

# Row 134 (G=44020)
$ws.Range("H134").Value = 7541.2666
$ws.Range("I134").Value = 5907.1904
$ws.Range("J134").Value = 11354.111
$ws.Range("K134").Value = 17721.5712
$ws.Range("L134").Value = 34062.333
$ws.Range("M134").Value = -15186.5712
$ws.Range("N134").Value = -39132.333

$ws = $wb.Worksheets.Item("CUL")
# Row 22 (G=4697)
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("K22").Value = 3000
$ws.Range("M22").Value = -2831

# Row 27 (G=4697)
$ws.Range("H27").Value = 1000
$ws.Range("I27").Value = 1000
$ws.Range("K27").Value = 3000
$ws.Range("M27").Value = -2898

# Row 81 (G=12843)
$ws.Range("H81").Value = 5644.1
$ws.Range("I81").Value = 5976.6
$ws.Range("J81").Value = 5311.6
$ws.Range("K81").Value = 17929.8
$ws.Range("L81").Value = 15934.8
$ws.Range("M81").Value = -16806.8
$ws.Range("N81").Value = -18180.8

# Row 84 (G=12843)
$ws.Range("H84").Value = 5644.1
$ws.Range("I84").Value = 5976.6
$ws.Range("J84").Value = 5311.6
$ws.Range("K84").Value = 53789.4
$ws.Range("L84").Value = 47804.4
$ws.Range("M84").Value = -48173.4
$ws.Range("N84").Value = -59036.4

# Row 114 (G=27865)
$ws.Range("H114").Value = 1410.6111
$ws.Range("I114").Value = 657.3333
$ws.Range("J114").Value = 1561.2667
$ws.Range("K114").Value = 1971.9999
$ws.Range("L114").Value = 4683.800099999999
$ws.Range("M114").Value = 1282.0001
$ws.Range("N114").Value = -11191.8001

# Row 117 (G=27870)
$ws.Range("H117").Value = 699.6667
$ws.Range("J117").Value = 799.6667
$ws.Range("L117").Value = 2399.0001
$ws.Range("N117").Value = -9283.000100000001

# Row 131 (G=36060)
$ws.Range("H131").Value = 17150.969
$ws.Range("J131").Value = 1406.0204
$ws.Range("L131").Value = 4218.0612
$ws.Range("N131").Value = -14298.0612

# Row 132 (G=43972)
$ws.Range("H132").Value = 5264542.5
$ws.Range("I132").Value = 1450.7778
$ws.Range("J132").Value = 10001325
$ws.Range("K132").Value = 13057.0002
$ws.Range("L132").Value = 90011925
$ws.Range("M132").Value = -10527.0002
$ws.Range("N132").Value = -90016985

# Row 134 (G=44074)
$ws.Range("H134").Value = 1303.3125
$ws.Range("I134").Value = 1303.3125
$ws.Range("K134").Value = 3909.9375
$ws.Range("M134").Value = 1160.0625

# Row 137 (G=44088)
$ws.Range("H137").Value = 3221.4546
$ws.Range("I137").Value = 2491.2856
$ws.Range("K137").Value = 7473.8568
$ws.Range("M137").Value = -2373.8568

$ws = $wb.Worksheets.Item("GSM")
# Row 107 (G=27802)
$ws.Range("H107").Value = 420.6154
$ws.Range("J107").Value = 616.2
$ws.Range("L107").Value = 616.2
$ws.Range("N107").Value = -4456.2

# Row 132 (G=44008)
$ws.Range("H132").Value = 3255
$ws.Range("I132").Value = 3088.16
$ws.Range("K132").Value = 9264.48
$ws.Range("M132").Value = -6734.48

$ws = $wb.Worksheets.Item("LTW")
# Row 20 (G=4308)
$ws.Range("H20").Value = 87930.97
$ws.Range("J20").Value = 90446.36
$ws.Range("L20").Value = 90446.36
$ws.Range("N20").Value = -90898.36

# Row 22 (G=5277)
$ws.Range("H22").Value = 3612.7144
$ws.Range("J22").Value = 6999.5
$ws.Range("L22").Value = 6999.5
$ws.Range("N22").Value = -7589.5

# Row 27 (G=5277)
$ws.Range("H27").Value = 3612.7144
$ws.Range("J27").Value = 6999.5
$ws.Range("L27").Value = 6999.5
$ws.Range("N27").Value = -7213.5

# Row 40 (G=36248)
$ws.Range("H40").Value = 2606.25
$ws.Range("I40").Value = 2218
$ws.Range("K40").Value = 2218
$ws.Range("M40").Value = -2082

# Row 63 (G=12006)
$ws.Range("H63").Value = 51271
$ws.Range("J63").Value = 25085
$ws.Range("L63").Value = 25085
$ws.Range("N63").Value = -26583

# Row 66 (G=12006)
$ws.Range("H66").Value = 51271
$ws.Range("J66").Value = 25085
$ws.Range("L66").Value = 75255
$ws.Range("N66").Value = -82743

# Row 132 (G=44058)
$ws.Range("H132").Value = 2593.12
$ws.Range("I132").Value = 2196.658
$ws.Range("K132").Value = 6589.974
$ws.Range("M132").Value = -4059.974

# Row 136 (G=44060)
$ws.Range("H136").Value = 3274.087
$ws.Range("I136").Value = 3107.4707
$ws.Range("K136").Value = 9322.4121
$ws.Range("M136").Value = -6772.4121

$ws = $wb.Worksheets.Item("WVR")
# Row 101 (G=18538)
$ws.Range("H101").Value = 28055.223
$ws.Range("J101").Value = 28055.223
$ws.Range("L101").Value = 28055.223
$ws.Range("N101").Value = -34545.223

# Row 109 (G=27161)
$ws.Range("H109").Value = 57443
$ws.Range("J109").Value = 57443
$ws.Range("L109").Value = 57443
$ws.Range("N109").Value = -60217

# Row 132 (G=44029)
$ws.Range("H132").Value = 19778.11
$ws.Range("I132").Value = 11664.139
$ws.Range("J132").Value = 35151.95
$ws.Range("K132").Value = 34992.417
$ws.Range("L132").Value = 105455.85
$ws.Range("M132").Value = -32462.417
$ws.Range("N132").Value = -110515.85

# Row 138 (G=42347)
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
